# Apply the OOXML diff to slide 1:
#  - Shape "Straight Arrow Connector 45" (id 46): shrink width slightly
#  - Shape "Straight Arrow Connector 50" (id 51): move/shrink slightly
#  - Shape "Rectangle: Rounded Corners 24" (id 25): reposition/resize, and
#    change the space after "rpId" to a comma.
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points (a
# 32-bit/Single float under the COM hood) while the OOXML stores English
# Metric Units (EMU, 12700 per point). Converting an EMU value naively to
# points and back can truncate to one EMU less than intended once it has
# been rounded through a Single, so EmuToPt nudges the point value by the
# smallest amount needed to land back exactly on the target EMU.

function EmuToPt {
    param([double]$Emu)
    $EmuPerPoint = 12700.0
    $base = $Emu / $EmuPerPoint
    for ($steps = 0; $steps -lt 5000; $steps++) {
        $candidate = $base + ($steps * 0.0000001)
        $asSingle = [float]$candidate
        $back = [math]::Floor([double]$asSingle * $EmuPerPoint)
        if ($back -eq $Emu) {
            return $candidate
        }
    }
    throw "EmuToPt: no suitable point value found for EMU=$Emu"
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Straight Arrow Connector 45 (index 35, id 46) ---
$conn1 = $s.Shapes.Item(35)
$conn1.Width = EmuToPt 2088000

# --- Straight Arrow Connector 50 (index 38, id 51) ---
$conn2 = $s.Shapes.Item(38)
$conn2.Left = EmuToPt 5983246
$conn2.Width = EmuToPt 3348000

# --- Rectangle: Rounded Corners 24 (index 50, id 25) ---
$rect = $s.Shapes.Item(50)

# Change the " " run right after "rpId" (and before the line break /
# "filtering issuers") into "," first. This shape auto-fits its height to
# its text (<a:spAutoFit/>), so resize the box afterwards to the diff's
# final values -- otherwise the autofit recalculation would clobber the
# explicit height we set below.
$tr = $rect.TextFrame.TextRange
$sub = $tr.Characters(36, 1)
$sub.Text = ","

$rect.Left = EmuToPt 4273630
$rect.Top = EmuToPt 3961833
$rect.Width = EmuToPt 1695663
$rect.Height = EmuToPt 2268000
